$wb = $excel.ActiveWorkbook

# Duplicate the "Czech" sheet to use as the template for the new "Swiss" sheet,
# placing the copy after the last existing sheet.
$czech = $wb.Worksheets.Item("Czech")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Copy($null, $lastSheet)

$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Update the Swiss-market specific values (order matches the new shared-string
# insertion order: Switzerland Market, P32AR-CH, P32DR-CH, NGC-3476/T2347).
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("A16").Value = "P32AR-CH"
$swiss.Range("A17").Value = "P32DR-CH"
$swiss.Range("B4").Value = "NGC-3476/T2347"

# The Swiss tab's columns B/D were manually resized (no longer auto "best fit")
# in the authored workbook; reproduce the saved widths as closely as possible.
$swiss.Columns.Item(2).ColumnWidth = 22.8307
$swiss.Columns.Item(4).ColumnWidth = 29.0534

# Match the saved selection/active-cell state from the authored workbook.
$swiss.Activate()
$swiss.Range("E13").Select()

$czech.Activate()
$czech.Range("A8:A17").Select()

$swiss.Activate()
